$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New lesson-plan content for Week06, filled in the same order the shared
# strings first appear in the target workbook.
$ws.Range("C3").Value = "Mi az az internet 2?"
$ws.Range("E3").Value = "Free wifi veszélyei."
$ws.Range("F4").Value = "Jelszókezelés, jelszó feltörések(databreaches)"
$ws.Range("D4").Value = "Email spam, phising, data theft, blackmail"
$ws.Range("D3").Value = "Mi az az email, internetes levelezés."
$ws.Range("F3").Value = "Jelszavak, kódolások működése, titkosítás alapjai"
$ws.Range("G3").Value = "Kvantumszámítógépek, titkosítás fontossága"
$ws.Range("C4").Value = "Beszélgetős óra - az előző órai dolgokról"
$ws.Range("E4").Value = "Beszélgetős óra - az előző órai dolgokról"
$ws.Range("G4").Value = "Beszélgetős óra - az előző órai dolgokról"
$ws.Range("G2").Value = "Week06"

# Column widths (bestFit-style) for the newly populated columns
$ws.Columns.Item(3).ColumnWidth = 36.333333333333336
$ws.Columns.Item(4).ColumnWidth = 42.333333333333336
$ws.Columns.Item(5).ColumnWidth = 36.333333333333336
$ws.Columns.Item(6).ColumnWidth = 44.5
$ws.Columns.Item(7).ColumnWidth = 40.166666666666664

# Remove now-unused L column cell (dimension shrinks from L13 to K13)
$ws.Range("L2").Clear()

# Update selection to match final state
$ws.Range("D11").Select()
